$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '91.147.55'
$ws.Range('E2').Value = '  +2.12%  '
$ws.Range('D3').Value = '3.141.41'
$ws.Range('E3').Value = '  +1.62%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.57'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '625.97'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.67%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.04'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +25.99%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.373'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.81%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').Value = '3.136.98'
$ws.Range('E10').Value = '  +1.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.756'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +23.38%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.198'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +8.92%  '
$ws.Range('E13').Value = '  +3.75%  '
$ws.Range('B14').Value = 'Toncoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.58'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.11%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '34.73'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.93%  '
$ws.Range('D16').Value = '90.976.57'
$ws.Range('E16').Value = '  +2.47%  '
$ws.Range('D17').Value = '3.712.42'
$ws.Range('E17').Value = '  +1.44%  '
$ws.Range('D18').Value = '3.134.71'
$ws.Range('E18').Value = '  +1.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.75'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +10.52%  '
$ws.Range('E20').Value = '  +1.74%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.21'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.62%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '444.77'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.60%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.90'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +7.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.22'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '6.25'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +13.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '89.13'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.61%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.34'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.26%  '
$ws.Range('D28').Value = '3.296.56'
$ws.Range('E28').Value = '  +1.27%  '
$ws.Range('E30').Value = '  -0.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.20'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +13.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '527.94'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.905'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -15.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.15'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +12.33%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.74'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.58%  '
$ws.Range('E36').Value = '  +13.71%  '
$ws.Range('E37').Value = '  +5.52%  '
$ws.Range('E38').Value = '  +5.04%  '
$ws.Range('E39').Value = '  +4.04%  '
$ws.Range('B40').Value = 'WhiteBITCoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '22.26'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.159'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +16.53%  '
$ws.Range('E42').Value = '  -0.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0822'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +17.85%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.404'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +11.13%  '
$ws.Range('B45').Value = 'USDe'
$ws.Range('C45').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.94'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '148.99'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.34%  '
$ws.Range('B48').Value = 'ImmutableX'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.33'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +8.63%  '
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '44.10'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.39'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +10.65%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '171.73'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.83%  '
